# Apply updated crypto price/volume data (rows 2-51) to match target snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.438.52'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '1.644.54'
$ws.Range("E3").Value = '  -1.23%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.25'
$ws.Range("E5").Value = '  -1.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.539'
$ws.Range("E6").Value = '  +4.84%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.27'
$ws.Range("E8").Value = '  -1.07%  '
$ws.Range("E9").Value = '  -2.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0610'
$ws.Range("E10").Value = '  -1.86%  '
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("D12").Value = '1.877.09'
$ws.Range("E12").Value = '  -1.28%  '
$ws.Range("D13").Value = '1.642.83'
$ws.Range("E13").Value = '  -1.42%  '
$ws.Range("E14").Value = '  -2.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.558'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.33'
$ws.Range("E16").Value = '  -2.76%  '
$ws.Range("D17").Value = '27.407.95'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.01'
$ws.Range("E18").Value = '  -7.64%  '
$ws.Range("E19").Value = '  -1.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.53'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("E22").Value = '  -3.75%  '
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.06'
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.115'
$ws.Range("E26").Value = '  +2.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.95'
$ws.Range("E27").Value = '  -2.92%  '
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.55'
$ws.Range("E29").Value = '  -5.19%  '
$ws.Range("E30").Value = '  -4.69%  '
$ws.Range("E31").Value = '  -3.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.28'
$ws.Range("E32").Value = '  -2.05%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '1.413.80'
$ws.Range("E34").Value = '  -4.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.57'
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("E36").Value = '  -0.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.883'
$ws.Range("E37").Value = '  -5.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.562'
$ws.Range("E38").Value = '  -1.99%  '
$ws.Range("E39").Value = '  -3.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.03'
$ws.Range("E40").Value = '  +0.59%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.48'
$ws.Range("E42").Value = '  -1.56%  '
$ws.Range("E43").Value = '  +1.19%  '
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.62'
$ws.Range("E45").Value = '  -7.11%  '
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").Value = '1.787.53'
$ws.Range("E47").Value = '  -1.15%  '
$ws.Range("E48").Value = '  -3.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '87.50'
$ws.Range("E49").Value = '  -1.96%  '
$ws.Range("D50").Value = '0.0₆0105'
$ws.Range("E50").Value = '  -3.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0988'
$ws.Range("E51").Value = '  -3.15%  '
